$d = $word.ActiveDocument

$old = "Estudar o material didático e as demais fontes de pesquisa proposta para o Curso."
$new = "Estudar o material didático e as demais fontes de pesquisa proposta para o Curso {cursos_FIC}."

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
